# Applies the cryptos.xlsx price/volume update described in the commit diff.
# Values are forced to Text via a leading apostrophe (matches the source
# workbook, where these cells are inlineStr / shared-string text, not numbers),
# then the cell Style is reset to "Normal" so no stray number-format style is
# introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''66.135.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +6.33%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''3.016.12'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +3.43%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.03%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''584.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  +2.67%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''162.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  +12.55%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  +0.03%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = '''XRP'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = '''https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").Value = '''0.521'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +4.16%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = '''LidoStakedEther'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = '''https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = '''3.012.56'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  +3.34%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''6.72'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -3.49%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.156'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  +4.40%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.457'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +5.82%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''0.0000255'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +6.31%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''34.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +6.20%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = '''  -0.65%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''66.118.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +6.35%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''3.517.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +3.43%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = '''  +5.37%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''3.016.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  +2.64%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''457.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +5.92%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''13.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +6.55%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '''  +5.50%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''7.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  +7.80%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''82.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  +4.55%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = '''  +14.77%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''12.38'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +2.88%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''10.55'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +3.78%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  -0.05%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''8.12'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +16.27%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''2.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +17.22%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''0.0000105'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -6.84%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  +4.10%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''27.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +6.21%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  +5.30%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -0.04%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''0.996'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  +4.30%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''5.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  +7.79%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''2.19'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  +15.50%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''3.01'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +2.91%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''50.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  +2.38%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.308'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +15.47%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = '''  +8.62%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''43.70'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +6.19%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  +3.71%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''398.21'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +14.12%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '''  +7.24%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''2.804.65'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +3.02%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''134.45'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  +0.66%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D50").Value = '''24.04'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  +11.49%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''0.108'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +4.66%  '
$ws.Range("E51").Style = "Normal"

Write-Host "Applied 90 cell updates across rows 2-51 (cryptos price refresh)."
